$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 252.0675830790046

$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 9.793184359356808

$ws.Range("B4").Value = 0.1554434735375247
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.7127328510149897
$ws.Range("E4").Value = 6.48142807727062
$ws.Range("G4").Value = 9.002840860721944

$ws.Range("B5").Value = 0.3464964993005633
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 157.8057217802531
$ws.Range("E5").Value = 246.9852506941017
$ws.Range("G5").Value = 406.7907054325541

$ws.Range("B6").Value = 0.06328177979961902
$ws.Range("C6").Value = 9.226618575922256
$ws.Range("D6").Value = 3.082599426703578
$ws.Range("E6").Value = 6.48142807727062
$ws.Range("G6").Value = 18.85392785969607
